# Adapt column header formatting to respective input file names.
#
# The sheet contains a diff table with columns suffixed "_old" (left side,
# format version FV2210) and "_new" (right side, format version FV2304).
# This renames those header cells to use the explicit format-version suffix,
# wraps the data range in a real Excel Table (ListObject), and freezes the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base (suffix-less) header names, in the order they appear for both the
# "_old"/FV2210 block (columns A-J) and the "_new"/FV2304 block (columns L-U).
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A..J (1..10) -> "<Name>_FV2210"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2210"
}

# Column K (11) is "diff" and stays unchanged.

# Columns L..U (12..21) -> "<Name>_FV2304"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2304"
}

# Turn the data range into a real table (Table1) with an autofilter, using
# the freshly renamed header row as the column names.
$dataRange = $ws.Range("A1:U88")
$tbl = $ws.ListObjects.Add(1, $dataRange, [System.Reflection.Missing]::Value, 1)
$tbl.TableStyle = ""

# Freeze the header row (split below row 1, keep the header visible while
# scrolling).
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
